$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the default (unstyled) cell format reference, used when forcing
# numeric-looking strings to remain text so the style index is not altered.
$defaultStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = "67.391.73"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "3.523.62"
$ws.Range("E3").Value = "  -0.51%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.56"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.72"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -1.50%  "

$ws.Range("D7").Value = "3.523.15"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("E11").Value = "  +2.61%  "

$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("E13").Value = "  -1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.09"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "4.117.29"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "3.522.81"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "67.395.71"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.08"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  -2.16%  "

$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("E23").Value = "  -2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.48"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  -0.64%  "

$ws.Range("E25").Value = "  +8.37%  "

$ws.Range("D26").Value = "3.662.76"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.34"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E32").Value = "  -6.73%  "

$ws.Range("E33").Value = "  +3.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.92"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").Value = "3.513.69"
$ws.Range("E36").Value = "  -0.81%  "

$ws.Range("E37").Value = "  -3.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.02"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.28"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +1.08%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("E42").Value = "  +4.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0880"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.44"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -3.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.882"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.62"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -3.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "44.89"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -2.17%  "

$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("E49").Value = "  +3.87%  "

$ws.Range("E50").Value = "  -1.36%  "

$ws.Range("E51").Value = "  -1.62%  "
